$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.117.99"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "2.929.39"
$ws.Range("E3").Value = "  +0.87%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'593.89"
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("D6").Value = "'146.25"
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.507"
$ws.Range("E8").Value = "  +0.36%  "
$ws.Range("D9").Value = "'6.90"
$ws.Range("E9").Value = "  +2.99%  "
$ws.Range("E10").Value = "  +0.68%  "
$ws.Range("E11").Value = "  -1.54%  "
$ws.Range("E12").Value = "  +1.16%  "
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("D15").Value = "3.411.84"
$ws.Range("D16").Value = "61.093.88"
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("E17").Value = "  -1.30%  "
$ws.Range("D18").Value = "2.929.42"
$ws.Range("E18").Value = "  +0.80%  "
$ws.Range("D19").Value = "'432.23"
$ws.Range("E19").Value = "  +1.37%  "
$ws.Range("E20").Value = "  -1.18%  "
$ws.Range("E21").Value = "  +2.19%  "
$ws.Range("D22").Value = "'7.10"
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("D23").Value = "'81.43"
$ws.Range("E23").Value = "  +1.22%  "
$ws.Range("E24").Value = "  +0.72%  "
$ws.Range("E25").Value = "  +0.96%  "
$ws.Range("D26").Value = "'12.01"
$ws.Range("E26").Value = "  +1.05%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("E28").Value = "  +5.83%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("D31").Value = "'7.08"
$ws.Range("E31").Value = "  -1.57%  "
$ws.Range("D32").Value = "'26.53"
$ws.Range("E32").Value = "  +0.22%  "
$ws.Range("E33").Value = "  +1.91%  "
$ws.Range("D34").Value = "0.0₃0854"
$ws.Range("E34").Value = "  +2.33%  "
$ws.Range("D35").Value = "'1.02"
$ws.Range("E35").Value = "  +1.38%  "
$ws.Range("E36").Value = "  -0.15%  "
$ws.Range("D37").Value = "'3.04"
$ws.Range("E37").Value = "  +3.82%  "
$ws.Range("E38").Value = "  +1.32%  "
$ws.Range("D39").Value = "'1.99"
$ws.Range("E39").Value = "  -1.51%  "
$ws.Range("D40").Value = "'8.61"
$ws.Range("E40").Value = "  -1.22%  "
$ws.Range("E41").Value = "  -0.86%  "
$ws.Range("D42").Value = "'40.34"
$ws.Range("E42").Value = "  -3.64%  "
$ws.Range("D43").Value = "'376.28"
$ws.Range("E43").Value = "  +0.96%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0347"
$ws.Range("E44").Value = "  +0.21%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.727.33"
$ws.Range("E45").Value = "  +2.54%  "
$ws.Range("D46").Value = "'130.76"
$ws.Range("E46").Value = "  -2.24%  "
$ws.Range("D48").Value = "'24.06"
$ws.Range("E48").Value = "  -3.26%  "
$ws.Range("E49").Value = "  +0.19%  "
$ws.Range("E50").Value = "  -2.50%  "
$ws.Range("E51").Value = "  +2.45%  "
